# Generate Report for Handoff
# Updates the handoff status text/timestamps and the now-narrower
# "Status" columns that result from the shorter status string.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refreshed "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" ---
$overview.Range("G2").Value = "2016-10-20 00:56:29"
$zhcn.Range("H2").Value     = "2016-10-20 00:56:18"
$dede.Range("H2").Value     = "2016-10-20 00:56:29"

# --- Narrow the Status columns to fit the shorter text ---
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth     = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth     = 16.3333333333333
